# Version EEE with web build for IOS
#
# Applies the "data" sheet update:
#  - every "Any" value in column C (the Grade column) becomes "Any Grade",
#    rendered in an explicit black font (new font + new cellXf get created
#    automatically by the engine, matching fonts/cellXfs growth in the diff)
#  - the AutoFilter / _FilterDatabase range moves from column D to column C
#  - the active selection on the sheet moves to E70

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Column C: "Any" -> "Any Grade" (with explicit black font colour) ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 70) { $lastRow = 70 }

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "Any") {
        $cell.Value = "Any Grade"
        $cell.Font.Color = 0
    }
}

# --- 2. Move the AutoFilter from D1:D70 to C1:C70 ---
$ws.AutoFilterMode = $false
[void]$ws.Range("C1:C70").AutoFilter()

# Keep the hidden _xlnm._FilterDatabase defined name in sync with the new
# filter range (Excel normally does this itself when the AutoFilter moves).
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=data!`$C`$1:`$C`$70"
    }
}

# --- 3. Update the saved selection on the sheet ---
[void]$ws.Range("E70").Select()
